$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = 1.62
$ws.Range("H5").Value = 3.8
$ws.Range("I5").Value = 5.5
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 7.5
$ws.Range("O5").Value = 1.4
$ws.Range("P5").Value = 2.75
$ws.Range("Q5").Value = 2.25
$ws.Range("R5").Value = 1.62
$ws.Range("AC5").Value = 7.5
$ws.Range("AE5").Value = 21
$ws.Range("AH5").Value = 12
$ws.Range("AW5").Value = 7
$ws.Range("Q6").Value = 2.08
$ws.Range("R6").Value = 1.73
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 3.25
$ws.Range("I8").Value = 2.25
$ws.Range("J8").Value = 3.55
$ws.Range("K8").Value = 2.05
$ws.Range("L8").Value = 2.87
$ws.Range("N8").Value = 8.1
$ws.Range("O8").Value = 1.28
$ws.Range("P8").Value = 3.05
$ws.Range("Q8").Value = 1.88
$ws.Range("R8").Value = 1.82
$ws.Range("S8").Value = 1.4
$ws.Range("T8").Value = 2.52
$ws.Range("U8").Value = 1.65
$ws.Range("V8").Value = 1.98
$ws.Range("W8").Value = 9.75
$ws.Range("Y8").Value = 10.5
$ws.Range("Z8").Value = 37
$ws.Range("AA8").Value = 25
$ws.Range("AB8").Value = 32
$ws.Range("AC8").Value = 10
$ws.Range("AD8").Value = 6.3
$ws.Range("AE8").Value = 13
$ws.Range("AF8").Value = 60
$ws.Range("AG8").Value = 450
$ws.Range("AH8").Value = 7.9
$ws.Range("AJ8").Value = 9
$ws.Range("AL8").Value = 18
$ws.Range("AN8").Value = 4.85
$ws.Range("AO8").Value = 16.5
$ws.Range("AQ8").Value = 75
$ws.Range("AT8").Value = 2.47
$ws.Range("AU8").Value = 6.9
$ws.Range("AV8").Value = 65
$ws.Range("AW8").Value = 4.15
$ws.Range("AX8").Value = 12
$ws.Range("AY8").Value = 20
$ws.Range("AZ8").Value = 50
$ws.Range("BA8").Value = 90
$ws.Range("BB8").Value = 250
$ws.Range("G9").Value = 2.15
$ws.Range("H9").Value = 3.1
$ws.Range("I9").Value = 3.5
$ws.Range("J9").Value = 3
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 4
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 7.5
$ws.Range("Q9").Value = 2.3
$ws.Range("R9").Value = 1.6
$ws.Range("S9").Value = 1.5
$ws.Range("T9").Value = 2.5
$ws.Range("U9").Value = 2
$ws.Range("V9").Value = 1.73
$ws.Range("Z9").Value = 21
$ws.Range("AA9").Value = 21
$ws.Range("AC9").Value = 7.5
$ws.Range("AE9").Value = 17
$ws.Range("AH9").Value = 9
$ws.Range("AO9").Value = 13
$ws.Range("AP9").Value = 26
$ws.Range("AT9").Value = 2.5
$ws.Range("AY9").Value = 34
$ws.Range("M10").Value = 1.03
$ws.Range("O10").Value = 1.3
$ws.Range("G11").Value = 4.33
$ws.Range("H11").Value = 3.3
$ws.Range("I11").Value = 1.83
$ws.Range("J11").Value = 5
$ws.Range("M11").Value = 1.03
$ws.Range("N11").Value = 8
$ws.Range("O11").Value = 1.34
$ws.Range("Q11").Value = 2.3
$ws.Range("R11").Value = 1.6
$ws.Range("S11").Value = 1.5
$ws.Range("T11").Value = 2.5
$ws.Range("U11").Value = 2.1
$ws.Range("V11").Value = 1.67
$ws.Range("Z11").Value = 51
$ws.Range("AB11").Value = 51
$ws.Range("AC11").Value = 7.5
$ws.Range("AI11").Value = 7.5
$ws.Range("AO11").Value = 26
$ws.Range("AP11").Value = 41
$ws.Range("AQ11").Value = 101
$ws.Range("AS11").Value = 351
$ws.Range("AT11").Value = 2.5
$ws.Range("AU11").Value = 9
$ws.Range("BA11").Value = 67
$ws.Range("G12").Value = 2
$ws.Range("H12").Value = 3.1
$ws.Range("I12").Value = 3.8
$ws.Range("M12").Value = 1.03
$ws.Range("O12").Value = 1.27
$ws.Range("P12").Value = 3.25
$ws.Range("Q12").Value = 2.08
$ws.Range("R12").Value = 1.73
$ws.Range("S12").Value = 1.44
$ws.Range("T12").Value = 2.63
$ws.Range("W12").Value = 7
$ws.Range("AD12").Value = 6
$ws.Range("AI12").Value = 21
$ws.Range("AT12").Value = 2.63
$ws.Range("BA12").Value = 101
$ws.Range("M13").Value = 1.02
$ws.Range("N13").Value = 11
$ws.Range("O13").Value = 1.22
$ws.Range("P13").Value = 3.5
$ws.Range("Q13").Value = 1.98
$ws.Range("R13").Value = 1.88
$ws.Range("G14").Value = 2.4
$ws.Range("H14").Value = 3.25
$ws.Range("I14").Value = 2.67
$ws.Range("L14").Value = 3.25
$ws.Range("U14").Value = 1.8
$ws.Range("V14").Value = 1.91
$ws.Range("AH14").Value = 8.25
$ws.Range("AJ14").Value = 10
$ws.Range("AK14").Value = 30
$ws.Range("AL14").Value = 23
$ws.Range("AO14").Value = 13
$ws.Range("AS14").Value = 300
$ws.Range("AW14").Value = 4.65
